# staging smoke test debugging
# Update the "build" column (A) for a subset of rows so that they report
# the latest build string ("33c392bc2b built at 2020-09-17 13:46") instead
# of the placeholder "test" value left over from staging.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newBuild = "33c392bc2b built at 2020-09-17 13:46`n"

$rows = @(6, 8, 9, 11, 12, 14, 15, 16, 17, 18, 20, 21, 23, 24)

foreach ($r in $rows) {
    $ws.Range("A$r").Value = $newBuild
    # The embedded newline makes Excel think the row needs extra height;
    # AutoFit keeps the row formatting untouched (matches the original file).
    $ws.Rows($r).AutoFit()
}
